$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 202
$wsOff.Range("C3").Value = 137
$wsOff.Range("D3").Value = 50
$wsOff.Range("E3").Value = 26

# Sheet "DEF" - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 186
$wsDef.Range("C3").Value = 117
$wsDef.Range("D3").Value = 54
$wsDef.Range("E3").Value = 31
$wsDef.Range("G3").Value = 2
